# Updates cryptos list: price (D) and 1h volume (E) columns,
# plus a row swap (Maker/VeChain) reflected as B/C/D/E content changes on rows 43-44.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.619.03"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "2.553.41"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.60"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.22"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("E12").Value = "  +8.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.50"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.96%  "

$ws.Range("D14").Value = "2.508.57"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.874"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.27%  "

$ws.Range("D17").Value = "42.689.61"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.29"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.78%  "

$ws.Range("D19").Value = "0.0₃0984"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.52"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.85"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.98%  "

$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "27.94"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.17%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.98"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.41%  "

$ws.Range("E29").Value = "  -0.47%  "

$ws.Range("E30").Value = "  +1.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.52"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.18"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.41%  "

$ws.Range("E33").Value = "  +1.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0800"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.29"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.21"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.40"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +14.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.10"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +33.03%  "

$ws.Range("E41").Value = "  +0.98%  "

$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.071.95"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0303"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.97%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.22"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.50%  "

$ws.Range("D48").Value = "2.801.89"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +8.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.06"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("E51").Value = "  +2.46%  "
